$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$log = ""

$cell = $tbl.Cell(1, 1)
$expected = "441×7=3087"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [1,1]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "959×5=4795"
}

$cell = $tbl.Cell(1, 2)
$expected = "742×9=6678"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [1,2]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "129×9=1161"
}

$cell = $tbl.Cell(1, 3)
$expected = "340×3=1020"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [1,3]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "459×8=3672"
}

$cell = $tbl.Cell(1, 4)
$expected = "405×5=2025"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [1,4]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "944×9=8496"
}

$cell = $tbl.Cell(1, 5)
$expected = "637×2=1274"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [1,5]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "565×4=2260"
}

$cell = $tbl.Cell(5, 1)
$expected = "532×6=3192"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [5,1]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "718×6=4308"
}

$cell = $tbl.Cell(5, 2)
$expected = "730×2=1460"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [5,2]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "437×3=1311"
}

$cell = $tbl.Cell(5, 3)
$expected = "527×7=3689"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [5,3]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "526×8=4208"
}

$cell = $tbl.Cell(5, 4)
$expected = "340×3=1020"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [5,4]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "327×5=1635"
}

$cell = $tbl.Cell(5, 5)
$expected = "936×8=7488"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [5,5]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "240×5=1200"
}

$cell = $tbl.Cell(10, 1)
$expected = "548×6=3288"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [10,1]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "231×5=1155"
}

$cell = $tbl.Cell(10, 2)
$expected = "698×7=4886"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [10,2]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "337×2=674"
}

$cell = $tbl.Cell(10, 3)
$expected = "625×2=1250"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [10,3]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "707×9=6363"
}

$cell = $tbl.Cell(10, 4)
$expected = "690×5=3450"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [10,4]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "651×4=2604"
}

$cell = $tbl.Cell(10, 5)
$expected = "492×4=1968"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [10,5]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "389×6=2334"
}

$cell = $tbl.Cell(15, 1)
$expected = "117×3=351"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [15,1]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "435×4=1740"
}

$cell = $tbl.Cell(15, 2)
$expected = "236×3=708"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [15,2]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "696×9=6264"
}

$cell = $tbl.Cell(15, 3)
$expected = "635×3=1905"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [15,3]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "337×4=1348"
}

$cell = $tbl.Cell(15, 4)
$expected = "292×5=1460"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [15,4]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "883×5=4415"
}

$cell = $tbl.Cell(15, 5)
$expected = "940×4=3760"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [15,5]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "238×3=714"
}

$cell = $tbl.Cell(20, 1)
$expected = "266×7=1862"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [20,1]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "855×4=3420"
}

$cell = $tbl.Cell(20, 2)
$expected = "297×5=1485"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [20,2]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "799×5=3995"
}

$cell = $tbl.Cell(20, 3)
$expected = "126×9=1134"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [20,3]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "692×7=4844"
}

$cell = $tbl.Cell(20, 4)
$expected = "258×8=2064"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [20,4]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "224×9=2016"
}

$cell = $tbl.Cell(20, 5)
$expected = "404×5=2020"
$actual = $cell.Range.Text
$actual = $actual.TrimEnd([char]13, [char]7)
if ($actual -ne $expected) {
    $log += "MISMATCH at [20,5]: expected $expected got $actual`n"
} else {
    $cell.Range.Text = "754×4=3016"
}

if ($log -eq "") { $log = "All replacements applied successfully." }
$log